$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/Volume(1h) updates are stored as plain text in the sheet (e.g. "301.05",
# "-0.91%"), so force the Text number format before assigning values to stop
# Excel from reinterpreting them as numbers/percentages, then clear the
# temporary formatting so the cell style matches the original (unstyled) cells.
$cells = @{
    "D2" = "301.05"
    "E2" = "-0.91%"
    "D3" = "31.43"
    "E3" = "-1.77%"
    "D4" = "5.159"
    "E4" = "-2.61%"
    "D5" = "0.07362"
    "E5" = "-1.40%"
    "D6" = "1.816"
    "E6" = "24.43%"
    "D7" = "7.835"
    "E7" = "0.54%"
    "D8" = "3.753"
    "E8" = "-1.49%"
    "D9" = "0.9297"
    "E9" = "0.76%"
    "D10" = "0.1708"
    "E10" = "0.71%"
    "D11" = "0.07099"
    "E11" = "-7.62%"
    "D12" = "0.08063"
    "E12" = "-0.20%"
    "D13" = "0.03035"
    "E13" = "-0.71%"
    "D14" = "0.09932"
    "E14" = "0.17%"
    "D15" = "0.001501"
    "E15" = "0.50%"
    "D16" = "0.006086"
    "E16" = "-5.90%"
    "D17" = "3.468"
    "E17" = "-0.38%"
    "E18" = "-0.33%"
    "E19" = "-1.73%"
    "D20" = "0.1317"
    "E20" = "-1.90%"
    "D21" = "4.579"
    "E21" = "1.78%"
    "D22" = "0.04645"
    "E22" = "1.85%"
    "D23" = "0.1581"
    "E23" = "-2.63%"
    "D24" = "0.001217"
    "E24" = "0.00%"
    "D25" = "0.004753"
    "E25" = "7.64%"
    "D26" = "0.0001297"
    "E26" = "-7.45%"
    "E27" = "7.33%"
    "D39" = "0.01722"
    "E39" = "-1.87%"
    "D40" = "0.04512"
    "E40" = "-0.65%"
    "D41" = "0.007104"
    "E41" = "-1.50%"
    "D42" = "0.1340"
    "E42" = "-0.32%"
    "D43" = "0.002186"
    "E43" = "-2.12%"
    "D44" = "0.01083"
    "E44" = "-14.93%"
    "D45" = "0.00006243"
    "E45" = "1.49%"
    "E46" = "-21.56%"
    "D47" = "1.845"
    "E47" = "160.46%"
}

$targetRange = $ws.Range("D2:E47")
$targetRange.NumberFormat = "@"
foreach ($addr in $cells.Keys) {
    $ws.Range($addr).Value = $cells[$addr]
}
$targetRange.ClearFormats()

